# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date serial in A1 by one day (2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in column D
$ws.Range("D33").Value = 214.542
$ws.Range("D34").Value = 237.684
$ws.Range("D35").Value = 274.92
$ws.Range("D39").Value = 293.364
$ws.Range("D40").Value = 441.09
